$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 20: date 2021-05-21, Developer1 value "4" (matches existing row19 text style)
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").Value2 = 44337
$ws.Range("B20").Value = "4"

# Row 21: date 2021-05-22, Developer1 value "6"
$ws.Range("A19").Copy()
$ws.Range("A21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A21").Value2 = 44338
$ws.Range("B21").Value = "6"

# Update the active selection to H23
$ws.Range("H23").Select()
$excel.CutCopyMode = $false
